$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5329.2
$ws.Range("I18").Value = 6782.5
$ws.Range("K18").Value = 6782.5
$ws.Range("M18").Value = -6498.5
$ws.Range("H32").Value = 1287.4706
$ws.Range("I32").Value = 772
$ws.Range("J32").Value = 1745.6666
$ws.Range("K32").Value = 772
$ws.Range("L32").Value = 1745.6666
$ws.Range("M32").Value = -446
$ws.Range("N32").Value = -2397.6666
$ws.Range("H40").Value = 4778283
$ws.Range("J40").Value = 6669096.5
$ws.Range("L40").Value = 6669096.5
$ws.Range("N40").Value = -6669446.5
$ws.Range("H69").Value = 11000
$ws.Range("I69").Value = 11000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 33000
$ws.Range("L69").Value = ""
$ws.Range("M69").Value = -32126
$ws.Range("N69").Value = 0
$ws.Range("H72").Value = 11000
$ws.Range("I72").Value = 11000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 99000
$ws.Range("L72").Value = ""
$ws.Range("M72").Value = -94632
$ws.Range("N72").Value = 0
$ws.Range("H88").Value = 23869096
$ws.Range("I88").Value = 111114344
$ws.Range("J88").Value = 74936.17999999999
$ws.Range("K88").Value = 111114344
$ws.Range("L88").Value = 74936.17999999999
$ws.Range("M88").Value = -111113938
$ws.Range("N88").Value = -75748.17999999999
$ws.Range("H91").Value = 23869096
$ws.Range("I91").Value = 111114344
$ws.Range("J91").Value = 74936.17999999999
$ws.Range("K91").Value = 111114344
$ws.Range("L91").Value = 74936.17999999999
$ws.Range("M91").Value = -111112940
$ws.Range("N91").Value = -77744.17999999999
$ws.Range("H112").Value = 5184.5454
$ws.Range("J112").Value = 5373.3096
$ws.Range("L112").Value = 16119.9288
$ws.Range("N112").Value = -18335.9288
$ws.Range("H132").Value = 3493.2415
$ws.Range("I132").Value = 3493
$ws.Range("K132").Value = 10479
$ws.Range("M132").Value = -7949
$ws.Range("H135").Value = 455145.28
$ws.Range("I135").Value = 476580.75
$ws.Range("K135").Value = 4289226.75
$ws.Range("M135").Value = -4286691.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = ""
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").Value = ""
$ws.Range("H61").Value = 5066.8096
$ws.Range("I61").Value = 2650.5
$ws.Range("K61").Value = 2650.5
$ws.Range("M61").Value = -2438.5
$ws.Range("H74").Value = 32224.096
$ws.Range("I74").Value = 38085.14
$ws.Range("K74").Value = 38085.14
$ws.Range("M74").Value = -37211.14
$ws.Range("H77").Value = 32224.096
$ws.Range("I77").Value = 38085.14
$ws.Range("K77").Value = 190425.7
$ws.Range("M77").Value = -186057.7
$ws.Range("H132").Value = 7732.892
$ws.Range("I132").Value = 6251
$ws.Range("K132").Value = 18753
$ws.Range("M132").Value = -16223
$ws.Range("H136").Value = 5066.8096
$ws.Range("I136").Value = 2650.5
$ws.Range("K136").Value = 7951.5
$ws.Range("M136").Value = -5401.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = ""
$ws.Range("H22").Value = 9325.182000000001
$ws.Range("I22").Value = 10237.7
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 10237.7
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -10064.7
$ws.Range("N22").Value = -546
$ws.Range("H26").Value = 29588.25
$ws.Range("I26").Value = 21468.666
$ws.Range("J26").Value = 53947
$ws.Range("K26").Value = 21468.666
$ws.Range("L26").Value = 53947
$ws.Range("M26").Value = -21176.666
$ws.Range("N26").Value = -54531
$ws.Range("H94").Value = 1449.4286
$ws.Range("I94").Value = 573.4286
$ws.Range("J94").Value = 3201.4285
$ws.Range("K94").Value = 573.4286
$ws.Range("L94").Value = 3201.4285
$ws.Range("M94").Value = -122.4286
$ws.Range("N94").Value = -4103.4285

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2699.5
$ws.Range("I2").Value = 2699.5
$ws.Range("K2").Value = 2699.5
$ws.Range("M2").Value = -2586.5
$ws.Range("H7").Value = 1419.1111
$ws.Range("J7").Value = 381.7143
$ws.Range("L7").Value = 381.7143
$ws.Range("N7").Value = -607.7143
$ws.Range("H16").Value = 5419
$ws.Range("I16").Value = 3650.75
$ws.Range("K16").Value = 3650.75
$ws.Range("M16").Value = -3363.75
$ws.Range("H21").Value = 3000
$ws.Range("I21").Value = 3000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 3000
$ws.Range("L21").Value = ""
$ws.Range("M21").Value = -2765
$ws.Range("N21").Value = 0
$ws.Range("H22").Value = 245.83333
$ws.Range("I22").Value = 240
$ws.Range("K22").Value = 240
$ws.Range("M22").Value = 110
$ws.Range("H31").Value = 6184.5146
$ws.Range("I31").Value = 2704.279
$ws.Range("K31").Value = 2704.279
$ws.Range("M31").Value = -2409.279
$ws.Range("H34").Value = 6184.5146
$ws.Range("I34").Value = 2704.279
$ws.Range("K34").Value = 2704.279
$ws.Range("M34").Value = -2502.279
$ws.Range("H57").Value = 30000
$ws.Range("J57").Value = 30000
$ws.Range("L57").Value = 30000
$ws.Range("N57").Value = -31120
$ws.Range("H86").Value = 5686537.5
$ws.Range("I86").Value = 10420827
$ws.Range("J86").Value = 5389.6
$ws.Range("K86").Value = 10420827
$ws.Range("L86").Value = 5389.6
$ws.Range("M86").Value = -10419704
$ws.Range("N86").Value = -7635.6
$ws.Range("H89").Value = 5686537.5
$ws.Range("I89").Value = 10420827
$ws.Range("J89").Value = 5389.6
$ws.Range("K89").Value = 52104135
$ws.Range("L89").Value = 26948
$ws.Range("M89").Value = -52098519
$ws.Range("N89").Value = -38180
$ws.Range("H107").Value = 1281.0454
$ws.Range("I107").Value = 849.6429000000001
$ws.Range("K107").Value = 849.6429000000001
$ws.Range("M107").Value = 1070.3571
$ws.Range("H113").Value = 5419
$ws.Range("I113").Value = 3650.75
$ws.Range("K113").Value = 3650.75
$ws.Range("M113").Value = -1480.75
$ws.Range("H132").Value = 6319.6113
$ws.Range("I132").Value = 2354.4167
$ws.Range("K132").Value = 7063.250100000001
$ws.Range("M132").Value = -4533.250100000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 11869.412
$ws.Range("I132").Value = 7533.5
$ws.Range("J132").Value = 15723.556
$ws.Range("K132").Value = 67801.5
$ws.Range("L132").Value = 141512.004
$ws.Range("M132").Value = -65271.5
$ws.Range("N132").Value = -146572.004
$ws.Range("H134").Value = 76215.21000000001
$ws.Range("I134").Value = 94273.91
$ws.Range("K134").Value = 282821.73
$ws.Range("M134").Value = -277751.73
$ws.Range("H137").Value = 118662.94
$ws.Range("I137").Value = 84114.336
$ws.Range("J137").Value = 201579.6
$ws.Range("K137").Value = 252343.008
$ws.Range("L137").Value = 604738.8
$ws.Range("M137").Value = -247243.008
$ws.Range("N137").Value = -614938.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4395.5
$ws.Range("I80").Value = 3996.8333
$ws.Range("J80").Value = 4993.5
$ws.Range("K80").Value = 3996.8333
$ws.Range("L80").Value = 4993.5
$ws.Range("M80").Value = -2998.8333
$ws.Range("N80").Value = -6989.5
$ws.Range("H83").Value = 4395.5
$ws.Range("I83").Value = 3996.8333
$ws.Range("J83").Value = 4993.5
$ws.Range("K83").Value = 19984.1665
$ws.Range("L83").Value = 24967.5
$ws.Range("M83").Value = -14992.1665
$ws.Range("N83").Value = -34951.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 11391.909
$ws.Range("I22").Value = 1179.8
$ws.Range("J22").Value = 19902
$ws.Range("K22").Value = 1179.8
$ws.Range("L22").Value = 19902
$ws.Range("M22").Value = -884.8
$ws.Range("N22").Value = -20492
$ws.Range("H27").Value = 11391.909
$ws.Range("I27").Value = 1179.8
$ws.Range("J27").Value = 19902
$ws.Range("K27").Value = 1179.8
$ws.Range("L27").Value = 19902
$ws.Range("M27").Value = -1072.8
$ws.Range("N27").Value = -20116
$ws.Range("H46").Value = 48531590
$ws.Range("I46").Value = 34482760
$ws.Range("J46").Value = 55556010
$ws.Range("K46").Value = 34482760
$ws.Range("L46").Value = 55556010
$ws.Range("M46").Value = -34482572
$ws.Range("N46").Value = -55556386
$ws.Range("H82").Value = 3259.25
$ws.Range("I82").Value = 3018.5
$ws.Range("K82").Value = 3018.5
$ws.Range("M82").Value = -2657.5
$ws.Range("H85").Value = 3259.25
$ws.Range("I85").Value = 3018.5
$ws.Range("K85").Value = 3018.5
$ws.Range("M85").Value = -1770.5
$ws.Range("H93").Value = 8223.333000000001
$ws.Range("I93").Value = 8500.5
$ws.Range("K93").Value = 8500.5
$ws.Range("M93").Value = -7252.5
$ws.Range("H132").Value = 15631582
$ws.Range("I132").Value = 23813586
$ws.Range("K132").Value = 71440758
$ws.Range("M132").Value = -71438228
$ws.Range("H136").Value = 7588.024
$ws.Range("I136").Value = 3682.2173
$ws.Range("J136").Value = 12316.105
$ws.Range("K136").Value = 11046.6519
$ws.Range("L136").Value = 36948.315
$ws.Range("M136").Value = -8496.651899999999
$ws.Range("N136").Value = -42048.315
